$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.965.82'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '3.825.03'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '704.41'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').Value = '171.55'
$ws.Range('E6').Value = '  -1.73%  '
$ws.Range('D7').Value = '3.824.07'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').Value = '7.36'
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').Value = '36.63'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '4.470.78'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = '3.821.25'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '70.920.60'
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').Value = '7.22'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '10.68'
$ws.Range('E21').Value = '  -4.53%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '493.26'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').Value = '0.736'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').Value = '85.18'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').Value = '10.60'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').Value = '12.10'
$ws.Range('E27').Value = '  -2.09%  '
$ws.Range('D28').Value = '2.08'
$ws.Range('E28').Value = '  -3.41%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('D31').Value = '7.41'
$ws.Range('E31').Value = '  -2.65%  '
$ws.Range('E32').Value = '  -4.35%  '
$ws.Range('D33').Value = '29.37'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('E34').Value = '  -3.84%  '
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('D36').Value = '3.787.40'
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('E39').Value = '  -3.07%  '
$ws.Range('D40').Value = '1.03'
$ws.Range('E40').Value = '  +3.16%  '
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('D42').Value = '3.30'
$ws.Range('E42').Value = '  -3.27%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '0.000314'
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('D46').Value = '163.26'
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').Value = '427.74'
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('D48').Value = '48.89'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('E51').Value = '  -2.34%  '
